$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New signature-style name under "Name:" (F2, merged with F3) ---
$ws.Range("F2").Value = "kihun han"
$ws.Range("F2").Font.Size = 28
$ws.Range("F2").Font.Name = "Kunstler Script"
$ws.Range("F2").Font.Family = 4
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").Borders.Item(8).LineStyle = 1

$ws.Range("F3").Font.Size = 28
$ws.Range("F3").Font.Name = "Kunstler Script"
$ws.Range("F3").Font.Family = 4
$ws.Range("F3").HorizontalAlignment = -4108

[void]$ws.Range("F2:F3").Merge()

# --- New work-log rows for 11.12 and 12.12 ---
[void]$ws.Range("B7:C7").Copy()
[void]$ws.Range("B21:C21").PasteSpecial(-4122)
[void]$ws.Range("B7:C7").Copy()
[void]$ws.Range("B22:C22").PasteSpecial(-4122)

$ws.Range("A21").Value = 11.12
$ws.Range("B21").Value = 0.375
$ws.Range("C21").Value = 0.66666666666666663
$ws.Range("E21").Value = "7hr"
$ws.Range("F21").Value = "implementation for project html file from function javascript and try to figure out function c ( truth table) more"

$ws.Range("A22").Value = 12.12
$ws.Range("B22").Value = 0.375
$ws.Range("C22").Value = 0.54166666666666663
$ws.Range("E22").Value = "4hr"
$ws.Range("F22").Value = "implementation for function c (truth table) "

[void]$ws.Range("F22").Select()

Write-Output "done"
